$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "GroupName2"
$ws.Range("A2").Value = "TestGroup2"
$ws.Range("A3").Value = "notes2"

$ws.Range("B2").Select()
